$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.830.86"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "'2.988.01"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'560.57"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").Value = "'137.48"
$ws.Range("E6").Value = "  +11.66%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +4.09%  "
$ws.Range("D9").Value = "'2.980.21"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").Value = "'0.132"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").Value = "'4.84"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  +7.75%  "
$ws.Range("D14").Value = "'33.73"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "'3.472.93"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "'6.97"
$ws.Range("E17").Value = "  +5.52%  "
$ws.Range("D18").Value = "'2.984.61"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("D19").Value = "'58.824.75"
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "'425.65"
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("D21").Value = "'13.54"
$ws.Range("E21").Value = "  +4.76%  "
$ws.Range("D22").Value = "'0.712"
$ws.Range("E22").Value = "  +5.63%  "
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("D24").Value = "'13.42"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("D25").Value = "'80.28"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'2.10"
$ws.Range("E28").Value = "  +8.04%  "
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "'7.74"
$ws.Range("E30").Value = "  +6.55%  "
$ws.Range("D31").Value = "'25.70"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").Value = "'6.08"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'0.0994"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "'0.985"
$ws.Range("E34").Value = "  +7.86%  "
$ws.Range("D35").Value = "'0.0₃0756"
$ws.Range("E35").Value = "  +21.37%  "
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("D38").Value = "'48.57"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "'8.72"
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("D40").Value = "'2.76"
$ws.Range("E40").Value = "  +15.26%  "
$ws.Range("D41").Value = "'397.73"
$ws.Range("E41").Value = "  +10.23%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").Value = "'2.724.25"
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("D45").Value = "'125.54"
$ws.Range("E45").Value = "  +4.59%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  +5.83%  "
$ws.Range("D48").Value = "'2.02"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "'32.11"
$ws.Range("E50").Value = "  +19.96%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'23.28"
$ws.Range("E51").Value = "  +2.05%  "
